$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows before row 7, pushing the old rows 7-10 down to 10-13
$ws.Range("A7:K9").EntireRow.Insert()

function Set-PlanetRow($r, $id, $radius, $mass, $temp, $k) {
    $ws.Cells.Item($r, 1).Value = $id
    $ws.Cells.Item($r, 2).Value = $radius
    $ws.Cells.Item($r, 3).Value = $mass
    $ws.Cells.Item($r, 4).Value = $temp
    $ws.Cells.Item($r, 5).Value = 10.35852619955
    $ws.Cells.Item($r, 6).Value = 531.944
    $ws.Cells.Item($r, 7).Value = 5590.19
    $ws.Cells.Item($r, 8).Value = 0
    $ws.Cells.Item($r, 9).Value = 0.103
    $ws.Cells.Item($r, 10).Value = "Unknown"
    $ws.Cells.Item($r, 11).Value = $k
}

# Fill the 3 new rows (7,8,9) with new planet data (ids 11, 13, 12)
Set-PlanetRow 7 11 1.1 1.05 300 0.974672675132751
Set-PlanetRow 8 13 1.1 1.05 300 0.974672675132751
Set-PlanetRow 9 12 1 1.1 300 0.972165942192078

# Old rows 7 and 8 (id 9, id 5) are now at rows 10 and 11 - keep them as-is.
# Old rows 9 and 10 (id 6, id 8) are now at rows 12 and 13 - delete them entirely.
$ws.Range("A12:K13").EntireRow.Delete()
